$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Insertion 1: new row for "Mailbox.getSelectedItemsAsync" before current row 141 ---
$ws.Rows.Item(141).Insert()
$ws.Range("A141").Value = "Office"
$ws.Range("B141").Value = "Mailbox"
$ws.Range("C141").Value = "getSelectedItemsAsync"
$ws.Range("D141").Value = 2
$ws.Range("E141").Value = "outlook-other-item-apis-get-message-properties"
$ws.Range("F141").Value = "getMessageProperties"
$ws.Range("E141").Style = "Normal"

# --- Insertion 2: new row for "SelectedItemDetails" interface before current row 285 (post-shift) ---
$ws.Rows.Item(285).Insert()
$ws.Range("A285").Value = "Office"
$ws.Range("B285").Value = "SelectedItemDetails"
$ws.Range("D285").Value = "interface"
$ws.Range("E285").Value = "outlook-other-item-apis-get-message-properties"
$ws.Range("F285").Value = "getMessageProperties"
$ws.Range("E285").Style = "Normal"

# --- Fix the previously mis-styled "getAsFileAsync" row, now shifted to row 233 ---
$ws.Range("C233:F233").Style = "Normal"
$ws.Range("D233").HorizontalAlignment = -4152

# --- Resize the table (ListObject) to cover the two newly inserted rows ---
$tbl.Resize($ws.Range("A1:F299"))

# --- Update the active selection to reflect where editing ended ---
$ws.Range("D285").Select()
